$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# The "INTRPT" / "RETINT" instructions in column M (rows 14-15) are being
# replaced by a richer interrupt documentation block in column Q
# (rows 2-6): *INTERRUPT, EXINT, NOINT, ENINT, SIV. Clear the now-stale
# column M entries for rows 14-17 and populate the new column Q cells.

$ws.Range("M14").Value = $null
$ws.Range("M15").Value = $null
$ws.Range("M16").Value = $null
$ws.Range("M17").Value = $null

$ws.Range("Q2").Value = "*INTERRUPT"
$ws.Range("Q3").Value = "EXINT"
$ws.Range("Q4").Value = "NOINT"
$ws.Range("Q5").Value = "ENINT"
$ws.Range("Q6").Value = "SIV"

$ws.Range("R4").Select()
